$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.003565333333333
$ws.Range("H2").Value = 3.010696
$ws.Range("I2").Value = 0.01537020623156171
$ws.Range("J2").Value = 0.01537020623156171
$ws.Range("M2").Value = 2.325008666666667
$ws.Range("N2").Value = 6.975026
$ws.Range("O2").Value = 0.05445297772988467
$ws.Range("P2").Value = 0.05445297772988466
$ws.Range("Q2").Value = 2.333298097566222
$ws.Range("R2").Value = 20.999682878096
$ws.Range("S2").Value = 0.0008369534976309645
$ws.Range("T2").Value = 0.0008369534976309642

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.003565333333333
$ws.Range("H3").Value = 3.010696
$ws.Range("I3").Value = 0.01537020623156171
$ws.Range("J3").Value = 0.01537020623156171
$ws.Range("O3").Value = 0.4529132218878514
$ws.Range("P3").Value = 0.4529132218878514
$ws.Range("Q3").Value = 19.40723176307645
$ws.Range("R3").Value = 174.665085867688
$ws.Range("S3").Value = 0.006961369625417347
$ws.Range("T3").Value = 0.006961369625417345

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.003565333333333
$ws.Range("H4").Value = 3.010696
$ws.Range("I4").Value = 0.01537020623156171
$ws.Range("J4").Value = 0.01537020623156171
$ws.Range("O4").Value = 0.492633800382264
$ws.Range("P4").Value = 0.492633800382264
$ws.Range("Q4").Value = 21.10924979953689
$ws.Range("R4").Value = 189.983248195832
$ws.Range("S4").Value = 0.007571883108513403
$ws.Range("T4").Value = 0.007571883108513401

$ws.Range("I5").Value = 0.9295661239816587
$ws.Range("J5").Value = 0.9295661239816587
$ws.Range("M5").Value = 2.325008666666667
$ws.Range("N5").Value = 6.975026
$ws.Range("O5").Value = 0.05445297772988467
$ws.Range("P5").Value = 0.05445297772988466
$ws.Range("Q5").Value = 141.1142333402531
$ws.Range("R5").Value = 1270.028100062278
$ws.Range("S5").Value = 0.05061764344762848
$ws.Range("T5").Value = 0.05061764344762847

$ws.Range("I6").Value = 0.9295661239816587
$ws.Range("J6").Value = 0.9295661239816587
$ws.Range("O6").Value = 0.4529132218878514
$ws.Range("P6").Value = 0.4529132218878514
$ws.Range("S6").Value = 0.421012788170335
$ws.Range("T6").Value = 0.4210127881703349

$ws.Range("I7").Value = 0.9295661239816587
$ws.Range("J7").Value = 0.9295661239816587
$ws.Range("O7").Value = 0.492633800382264
$ws.Range("P7").Value = 0.492633800382264
$ws.Range("S7").Value = 0.4579356923636954
$ws.Range("T7").Value = 0.4579356923636953

$ws.Range("I8").Value = 0.05506366978677964
$ws.Range("J8").Value = 0.05506366978677963
$ws.Range("M8").Value = 2.325008666666667
$ws.Range("N8").Value = 6.975026
$ws.Range("O8").Value = 0.05445297772988467
$ws.Range("P8").Value = 0.05445297772988466
$ws.Range("Q8").Value = 8.359026158977779
$ws.Range("R8").Value = 75.2312354308
$ws.Range("S8").Value = 0.002998380784625235
$ws.Range("T8").Value = 0.002998380784625234

$ws.Range("I9").Value = 0.05506366978677964
$ws.Range("J9").Value = 0.05506366978677963
$ws.Range("O9").Value = 0.4529132218878514
$ws.Range("P9").Value = 0.4529132218878514
$ws.Range("S9").Value = 0.02493906409209911
$ws.Range("T9").Value = 0.0249390640920991

$ws.Range("I10").Value = 0.05506366978677964
$ws.Range("J10").Value = 0.05506366978677963
$ws.Range("O10").Value = 0.492633800382264
$ws.Range("P10").Value = 0.492633800382264
$ws.Range("S10").Value = 0.0271262249100553
$ws.Range("T10").Value = 0.0271262249100553

